# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.069.79'
$ws.Range('E2').Value = '  +0.74%  '
# Row 3: Ethereum
$ws.Range('D3').Value = '1.890.52'
$ws.Range('E3').Value = '  +1.55%  '
# Row 4: TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
# Row 6: USDC
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.03%  '
# Row 7: XRP
$ws.Range('E7').Value = '  +2.76%  '
# Row 8: Cardano
$ws.Range('E8').Value = '  +3.04%  '
# Row 9: Dogecoin
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07214'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.69%  '
# Row 10: Solana
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.50%  '
# Row 11: Polygon
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9021'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.16%  '
# Row 12: TRON
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07636'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.76%  '
# Row 13: WrappedEther
$ws.Range('D13').Value = '1.870.87'
$ws.Range('E13').Value = '  +0.48%  '
# Row 14: Litecoin
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.28%  '
# Row 15: Polkadot
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.241'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.27%  '
# Row 16: BinanceUSD
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9997'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.00%  '
# Row 17: ShibaInu
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008506'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.10%  '
# Row 18: Avalanche
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.46'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.79%  '
# Row 19: Dai
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9993'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.03%  '
# Row 20: WrappedBTC
$ws.Range('D20').Value = '27.110.87'
$ws.Range('E20').Value = '  +0.70%  '
# Row 21: Uniswap
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.053'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.53%  '
# Row 22: WrappedliquidstakedEther2.0
$ws.Range('D22').Value = '2.117.33'
$ws.Range('E22').Value = '  +0.51%  '
# Row 23: Cosmos
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.21%  '
# Row 24: Chainlink
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.391'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.16%  '
# Row 25: LidoDAOToken
$ws.Range('E25').Value = '  +10.24%  '
# Row 26: Monero
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.49%  '
# Row 27: Toncoin
$ws.Range('E27').Value = '  -2.38%  '
# Row 28: EthereumClassic
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.17%  '
# Row 29: BitcoinCash
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.82%  '
# Row 30: Filecoin
$ws.Range('E30').Value = '  +5.23%  '
# Row 31: InternetComputer(DFINITY)
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.798'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.97%  '
# Row 32: Stellar
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09198'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.24%  '
# Row 33: Hedera
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05046'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.88%  '
# Row 34: ARBITRUM
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.243'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.90%  '
# Row 35: ImmutableX
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7663'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.15%  '
# Row 36: HuobiToken
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.976'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.80%  '
# Row 37: MXToken
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.275'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.75%  '
# Row 38: RenderToken
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.600'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.09%  '
# Row 39: TheSandbox
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5603'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.53%  '
# Row 40: VeChain
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01991'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.52%  '
# Row 41: TrustWalletToken
$ws.Range('E41').Value = '  +0.08%  '
# Row 42: Aptos
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.062'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.82%  '
# Row 43: FraxShare
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.626'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.06%  '
# Row 44: Quant
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '118.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.45%  '
# Row 45: Algorand
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1508'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.59%  '
# Row 46: Decentraland
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4827'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.83%  '
# Row 47: PaxDollar
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.06%  '
# Row 48: EnergySwap
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.85%  '
# Row 49: NEARProtocol
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.598'
$ws.Range('D49').Style = 'Normal'
# Row 50: Elrond
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.85%  '
# Row 51: Aave
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.54%  '
